$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.357.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4702"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2906"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08033"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.877.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.150"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6871"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.360.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007640"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.123.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.308"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.218"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.266"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.964"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.369"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09904"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.377"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.467"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.082"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04718"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7100"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01883"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.636"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.336"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.966"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4180"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8436"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.279"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.116"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "934.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05688"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.91%  "
